# Append: 2025-12-11 12:54 JST
# Update the "取得日時" (retrieved at) timestamp in column A for all
# existing data rows (rows 2-12) on the "ランサーズ" sheet from
# "2025-12-11 12:40:37" to "2025-12-11 12:54:10".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-11 12:40:37"
$newTimestamp = "2025-12-11 12:54:10"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
